$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their existing text format (they are stored as text strings,
# e.g. "130.20" or "0.550", which must not be reinterpreted as numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.678.62'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.286.69'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.42%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '503.75'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.20'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.39%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.299.51'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0961'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.15%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.09%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +4.42%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.61%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +6.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.691.18'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '54.660.41'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.29%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.288.60'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.32'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.17'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '306.72'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.61%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.996'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.80'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.996'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.48%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.41'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '171.24'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0709'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +5.11%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.26%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.60%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.75%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.995'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.939'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.07%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.04%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.42%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.74%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.08'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +8.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.41'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '125.44'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.55%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0899'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '246.84'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +5.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.550'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.80%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.81%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.88%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.78%  '
